# Streamlit Approval System: normalize COST_CENTER / LEDGER_NAME / LEDGER_UNDER /
# TO / BY columns (AK:AO) on rows 2-19 to the literal text "0" instead of the
# numeric 0 (or blank) they held before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("AK", "AL", "AM", "AN", "AO")

for ($row = 2; $row -le 19; $row++) {
    foreach ($col in $columns) {
        $cell = $ws.Range("$col$row")
        # Force text interpretation so "0" is stored as a string, not a number...
        $cell.NumberFormat = "@"
        $cell.Value = "0"
        # ...then drop back to the default style so no formatting residue is left.
        $cell.Style = "Normal"
    }
}
